$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16
$ws.Range("Q16").Value = 610384
$ws.Range("R16").Value = 7121170

# Row 17
$ws.Range("Q17").Value = 610054
$ws.Range("R17").Value = 7121273

# Row 18
$ws.Range("A18").Value = 111814591
$ws.Range("AB18").Value = "17:50"
$ws.Range("B18").Value = 77515
$ws.Range("E18").Value = 6425
$ws.Range("F18").Value = "Garnlav"
$ws.Range("G18").Value = "Alectoria sarmentosa"
$ws.Range("H18").Value = "(Ach.) Ach."
$ws.Range("Q18").Value = 610012
$ws.Range("R18").Value = 7121464
$ws.Range("Z18").Value = "17:50"

# Row 19
$ws.Range("Q19").Value = 610155
$ws.Range("R19").Value = 7121461

# Row 20
$ws.Range("A20").Value = 111814688
$ws.Range("AB20").Value = "17:55"
$ws.Range("B20").Value = 90087
$ws.Range("D20").Value = "LC"
$ws.Range("E20").Value = 3298
$ws.Range("F20").Value = "Trådticka"
$ws.Range("G20").Value = "Climacocystis borealis"
$ws.Range("H20").Value = "(Fr.) Kotl. & Pouzar"
$ws.Range("Q20").Value = 610011
$ws.Range("R20").Value = 7121476
$ws.Range("Z20").Value = "17:55"

# Row 21
$ws.Range("Q21").Value = 610155
$ws.Range("R21").Value = 7121460

# Row 22
$ws.Range("A22").Value = 111815024
$ws.Range("AB22").Value = "18:12"
$ws.Range("B22").Value = 56414
$ws.Range("D22").Value = "NT"
$ws.Range("E22").Value = 100049
$ws.Range("F22").Value = "Spillkråka"
$ws.Range("G22").Value = "Dryocopus martius"
$ws.Range("H22").Value = "(Linnaeus, 1758)"
$ws.Range("Q22").Value = 609922
$ws.Range("R22").Value = 7121488
$ws.Range("Z22").Value = "18:12"

# Row 23
$ws.Range("Q23").Value = 610384
$ws.Range("R23").Value = 7121170

# Row 24
$ws.Range("AB24").ClearContents()
$ws.Range("Q24").Value = 610409
$ws.Range("R24").Value = 7121114
$ws.Range("Z24").ClearContents()

# Row 25
$ws.Range("A25").Value = 112013699
$ws.Range("AB25").Value = "19:40"
$ws.Range("B25").Value = 77515
$ws.Range("E25").Value = 6425
$ws.Range("F25").Value = "Garnlav"
$ws.Range("G25").Value = "Alectoria sarmentosa"
$ws.Range("H25").Value = "(Ach.) Ach."
$ws.Range("Q25").Value = 610068
$ws.Range("R25").Value = 7121408
$ws.Range("Z25").Value = "19:40"

# Row 26
$ws.Range("A26").Value = 112013691
$ws.Range("AB26").Value = "19:29"
$ws.Range("B26").Value = 88489
$ws.Range("E26").Value = 1962
$ws.Range("F26").Value = "Vaddporing"
$ws.Range("G26").Value = "Anomoporia kamtschatica"
$ws.Range("H26").Value = "(Parmasto) Bondartseva"
$ws.Range("Q26").Value = 610134
$ws.Range("R26").Value = 7121461
$ws.Range("Z26").Value = "19:29"

# Row 27
$ws.Range("A27").Value = 112013704
$ws.Range("AB27").Value = "19:49"
$ws.Range("B27").Value = 81248
$ws.Range("E27").Value = 1312
$ws.Range("F27").Value = "Gammelgransskål"
$ws.Range("G27").Value = "Pseudographis pinicola"
$ws.Range("H27").Value = "(Nyl.) Rehm"
$ws.Range("Q27").Value = 610094
$ws.Range("R27").Value = 7121455
$ws.Range("Z27").Value = "19:49"

# Row 28
$ws.Range("A28").Value = 112013703
$ws.Range("AB28").Value = "19:28"
$ws.Range("B28").Value = 77515
$ws.Range("E28").Value = 6425
$ws.Range("F28").Value = "Garnlav"
$ws.Range("G28").Value = "Alectoria sarmentosa"
$ws.Range("H28").Value = "(Ach.) Ach."
$ws.Range("Q28").Value = 610144
$ws.Range("R28").Value = 7121461
$ws.Range("Z28").Value = "19:28"

# Row 29
$ws.Range("A29").Value = 112013700
$ws.Range("AB29").Value = "19:35"
$ws.Range("Q29").Value = 610102
$ws.Range("R29").Value = 7121416
$ws.Range("Z29").Value = "19:35"

# Row 30
$ws.Range("A30").Value = 112013698
$ws.Range("AB30").Value = "19:49"
$ws.Range("B30").Value = 77515
$ws.Range("E30").Value = 6425
$ws.Range("F30").Value = "Garnlav"
$ws.Range("G30").Value = "Alectoria sarmentosa"
$ws.Range("H30").Value = "(Ach.) Ach."
$ws.Range("Q30").Value = 610094
$ws.Range("R30").Value = 7121456
$ws.Range("Z30").Value = "19:49"

# Row 31
$ws.Range("A31").Value = 112013696
$ws.Range("B31").Value = 86961
$ws.Range("E31").Value = 4962
$ws.Range("F31").Value = "Mjölsvärting"
$ws.Range("G31").Value = "Lyophyllum semitale"
$ws.Range("H31").Value = "(Fr. : Fr.) Kühner"
$ws.Range("Q31").Value = 610070
$ws.Range("R31").Value = 7121402

# Row 32
$ws.Range("A32").Value = 112013697
$ws.Range("AB32").Value = "19:35"
$ws.Range("B32").Value = 89423
$ws.Range("E32").Value = 5432
$ws.Range("F32").Value = "Granticka"
$ws.Range("G32").Value = "Porodaedalea chrysoloma"
$ws.Range("H32").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("Q32").Value = 610102
$ws.Range("R32").Value = 7121413
$ws.Range("Z32").Value = "19:35"

# Row 33
$ws.Range("A33").Value = 112013690
$ws.Range("AB33").Value = "19:43"
$ws.Range("B33").Value = 88489
$ws.Range("E33").Value = 1962
$ws.Range("F33").Value = "Vaddporing"
$ws.Range("G33").Value = "Anomoporia kamtschatica"
$ws.Range("H33").Value = "(Parmasto) Bondartseva"
$ws.Range("Q33").Value = 610052
$ws.Range("R33").Value = 7121425
$ws.Range("Z33").Value = "19:43"
